$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("Archer Heights")
$ws.Range("K6").Value = 33
$ws.Range("K7").Value = 87

$ws = $wb.Worksheets.Item("Ashburn")
$ws.Range("K4").Value = 11
$ws.Range("K7").Value = 170

$ws = $wb.Worksheets.Item("Auburn Gresham")
$ws.Range("K3").Value = 233
$ws.Range("K6").Value = 199
$ws.Range("K7").Value = 728

$ws = $wb.Worksheets.Item("Austin")
$ws.Range("K2").Value = 436
$ws.Range("K3").Value = 474
$ws.Range("K6").Value = 515
$ws.Range("K7").Value = 1559

$ws = $wb.Worksheets.Item("Avalon Park")
$ws.Range("K3").Value = 38
$ws.Range("K7").Value = 107

$ws = $wb.Worksheets.Item("Boystown")
$ws.Range("K4").Value = 7
$ws.Range("K6").Value = 34

$ws = $wb.Worksheets.Item("By Neighborhood")
$ws.Range("K4").Value = 87
$ws.Range("K6").Value = 170
$ws.Range("K7").Value = 728
$ws.Range("K8").Value = 1559
$ws.Range("K9").Value = 107
$ws.Range("K13").Value = 34
$ws.Range("K18").Value = 160
$ws.Range("K19").Value = 704
$ws.Range("K20").Value = 576
$ws.Range("K27").Value = 224
$ws.Range("K28").Value = 9
$ws.Range("K29").Value = 1298
$ws.Range("K31").Value = 265
$ws.Range("K33").Value = 1022
$ws.Range("K37").Value = 804
$ws.Range("K41").Value = 167
$ws.Range("K42").Value = 881
$ws.Range("K48").Value = 309
$ws.Range("K49").Value = 132
$ws.Range("K51").Value = 300
$ws.Range("K52").Value = 628
$ws.Range("K54").Value = 462
$ws.Range("K55").Value = 256
$ws.Range("K65").Value = 551
$ws.Range("K67").Value = 921
$ws.Range("K73").Value = 212
$ws.Range("K74").Value = 25
$ws.Range("K75").Value = 73
$ws.Range("K79").Value = 587
$ws.Range("K83").Value = 507
$ws.Range("K85").Value = 1097
$ws.Range("K89").Value = 356
$ws.Range("K90").Value = 225
$ws.Range("K93").Value = 89
$ws.Range("K94").Value = 318
$ws.Range("K95").Value = 392
$ws.Range("K97").Value = 187
$ws.Range("K98").Value = 120
$ws.Range("K99").Value = 406
$ws.Range("K101").Value = 23798

$ws = $wb.Worksheets.Item("Calumet Heights")
$ws.Range("K2").Value = 44
$ws.Range("K3").Value = 54
$ws.Range("K7").Value = 160

$ws = $wb.Worksheets.Item("Chatham")
$ws.Range("K2").Value = 206
$ws.Range("K3").Value = 210
$ws.Range("K5").Value = 22
$ws.Range("K6").Value = 233
$ws.Range("K7").Value = 704

$ws = $wb.Worksheets.Item("Chicago Lawn")
$ws.Range("K3").Value = 187
$ws.Range("K6").Value = 157
$ws.Range("K7").Value = 576

$ws = $wb.Worksheets.Item("Citywide Totals")
$ws.Range("K2").Value = 6877
$ws.Range("K3").Value = 7125
$ws.Range("K4").Value = 1468
$ws.Range("K5").Value = 507
$ws.Range("K6").Value = 7821
$ws.Range("K7").Value = 23798

$ws = $wb.Worksheets.Item("Edgewater")
$ws.Range("K6").Value = 81
$ws.Range("K7").Value = 224

$ws = $wb.Worksheets.Item("Edison Park")
$ws.Range("K6").Value = 2
$ws.Range("K7").Value = 9

$ws = $wb.Worksheets.Item("Englewood")
$ws.Range("K2").Value = 366
$ws.Range("K3").Value = 463
$ws.Range("K6").Value = 378
$ws.Range("K7").Value = 1298

$ws = $wb.Worksheets.Item("Gage Park")
$ws.Range("K3").Value = 68
$ws.Range("K6").Value = 94
$ws.Range("K7").Value = 265

$ws = $wb.Worksheets.Item("Garfield Park")
$ws.Range("K3").Value = 363
$ws.Range("K7").Value = 1022

$ws = $wb.Worksheets.Item("Grand Crossing")
$ws.Range("K2").Value = 232
$ws.Range("K6").Value = 239
$ws.Range("K7").Value = 804

$ws = $wb.Worksheets.Item("Hermosa")
$ws.Range("K6").Value = 67
$ws.Range("K7").Value = 167

$ws = $wb.Worksheets.Item("Humboldt Park")
$ws.Range("K2").Value = 239
$ws.Range("K6").Value = 329
$ws.Range("K7").Value = 881

$ws = $wb.Worksheets.Item("Lake View")
$ws.Range("K6").Value = 146
$ws.Range("K7").Value = 309

$ws = $wb.Worksheets.Item("Lincoln Park")
$ws.Range("K6").Value = 66
$ws.Range("K7").Value = 132

$ws = $wb.Worksheets.Item("Little Italy, UIC")
$ws.Range("K2").Value = 82
$ws.Range("K3").Value = 81
$ws.Range("K6").Value = 100
$ws.Range("K7").Value = 300

$ws = $wb.Worksheets.Item("Little Village")
$ws.Range("K2").Value = 168
$ws.Range("K6").Value = 230
$ws.Range("K7").Value = 628

$ws = $wb.Worksheets.Item("Loop")
$ws.Range("K6").Value = 251
$ws.Range("K7").Value = 462

$ws = $wb.Worksheets.Item("Lower West Side")
$ws.Range("K6").Value = 89
$ws.Range("K7").Value = 256

$ws = $wb.Worksheets.Item("New City")
$ws.Range("K2").Value = 181
$ws.Range("K7").Value = 551

$ws = $wb.Worksheets.Item("North Lawndale")
$ws.Range("K2").Value = 249
$ws.Range("K3").Value = 335
$ws.Range("K5").Value = 24
$ws.Range("K7").Value = 921

$ws = $wb.Worksheets.Item("Portage Park")
$ws.Range("K6").Value = 72
$ws.Range("K7").Value = 212

$ws = $wb.Worksheets.Item("Printers Row")
$ws.Range("K6").Value = 15
$ws.Range("K7").Value = 25

$ws = $wb.Worksheets.Item("Pullman")
$ws.Range("K2").Value = 24
$ws.Range("K7").Value = 73

$ws = $wb.Worksheets.Item("Roseland")
$ws.Range("K4").Value = 37
$ws.Range("K7").Value = 587

$ws = $wb.Worksheets.Item("South Chicago")
$ws.Range("K6").Value = 119
$ws.Range("K7").Value = 507

$ws = $wb.Worksheets.Item("South Shore")
$ws.Range("K2").Value = 360
$ws.Range("K3").Value = 383
$ws.Range("K6").Value = 267
$ws.Range("K7").Value = 1097

$ws = $wb.Worksheets.Item("Uptown")
$ws.Range("K6").Value = 105
$ws.Range("K7").Value = 356

$ws = $wb.Worksheets.Item("Washington Heights")
$ws.Range("K6").Value = 58
$ws.Range("K7").Value = 225

$ws = $wb.Worksheets.Item("West Lawn")
$ws.Range("K3").Value = 20
$ws.Range("K7").Value = 89

$ws = $wb.Worksheets.Item("West Loop")
$ws.Range("K6").Value = 146
$ws.Range("K7").Value = 318

$ws = $wb.Worksheets.Item("West Pullman")
$ws.Range("K2").Value = 135
$ws.Range("K7").Value = 392

$ws = $wb.Worksheets.Item("West Town")
$ws.Range("K6").Value = 100
$ws.Range("K7").Value = 187

$ws = $wb.Worksheets.Item("Wicker Park")
$ws.Range("K4").Value = 7
$ws.Range("K7").Value = 120

$ws = $wb.Worksheets.Item("Woodlawn")
$ws.Range("K2").Value = 108
$ws.Range("K3").Value = 169
$ws.Range("K7").Value = 406
